# Update the arithmetic answers in the single table of the document.
# The table has 20 rows x 5 columns; only every 4th row (1, 5, 9, 13, 17)
# actually holds the answer text, 25 cells total, matching the diff order.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "50÷7=7, 1"
$t.Cell(1, 2).Range.Text = "10÷6=1, 4"
$t.Cell(1, 3).Range.Text = "11÷2=5, 1"
$t.Cell(1, 4).Range.Text = "36÷5=7, 1"
$t.Cell(1, 5).Range.Text = "87÷7=12, 3"

$t.Cell(5, 1).Range.Text = "63÷6=10, 3"
$t.Cell(5, 2).Range.Text = "77÷4=19, 1"
$t.Cell(5, 3).Range.Text = "74÷4=18, 2"
$t.Cell(5, 4).Range.Text = "73÷2=36, 1"
$t.Cell(5, 5).Range.Text = "30÷3=10, 0"

$t.Cell(9, 1).Range.Text = "45÷2=22, 1"
$t.Cell(9, 2).Range.Text = "91÷8=11, 3"
$t.Cell(9, 3).Range.Text = "93÷5=18, 3"
$t.Cell(9, 4).Range.Text = "58÷2=29, 0"
$t.Cell(9, 5).Range.Text = "72÷3=24, 0"

$t.Cell(13, 1).Range.Text = "65÷7=9, 2"
$t.Cell(13, 2).Range.Text = "39÷3=13, 0"
$t.Cell(13, 3).Range.Text = "51÷4=12, 3"
$t.Cell(13, 4).Range.Text = "88÷5=17, 3"
$t.Cell(13, 5).Range.Text = "66÷8=8, 2"

$t.Cell(17, 1).Range.Text = "13÷8=1, 5"
$t.Cell(17, 2).Range.Text = "96÷2=48, 0"
$t.Cell(17, 3).Range.Text = "95÷2=47, 1"
$t.Cell(17, 4).Range.Text = "62÷7=8, 6"
$t.Cell(17, 5).Range.Text = "18÷4=4, 2"

Write-Output "Done updating table answers."
